# history data and live data donut
# Rename the "Color" header (E1) to lowercase "color" and drop the page
# margins to zero, matching the re-saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E1 held the shared string "Color" -> change it to "color"
$ws.Range("E1").Value = "color"

# Make E1 the active selection (matches the re-saved sheetView state)
$ws.Range("E1").Select() | Out-Null

# Zero out the page margins on the sheet
$ws.PageSetup.LeftMargin = 0
$ws.PageSetup.RightMargin = 0
$ws.PageSetup.TopMargin = 0
$ws.PageSetup.BottomMargin = 0
$ws.PageSetup.HeaderMargin = 0
$ws.PageSetup.FooterMargin = 0
